$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GIN")
$ws.Range("A58:A61").EntireRow.Insert()
$ws.Range("A66:A72").EntireRow.Insert()
$ws.Range("F58").Copy($ws.Range("F66:F72"))
Write-Output "ok"
